# Daily attendance processing - 2025-10-13 09:47:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "G2";  Value = "System, backup@backdoor.com, system" },
    @{ Cell = "G3";  Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G4";  Value = "System, backup@backdoor.com" },
    @{ Cell = "G5";  Value = "System, backup@backdoor.com" },
    @{ Cell = "G6";  Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G10"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "L10"; Value = "70.2%" },
    @{ Cell = "G11"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G12"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G13"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G14"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G15"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "S15"; Value = "69.4%" },
    @{ Cell = "S16"; Value = "67.6%" },
    @{ Cell = "G17"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "S17"; Value = "60.4%" },
    @{ Cell = "G18"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G19"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "H19"; Value = "35/53" },
    @{ Cell = "G29"; Value = "System, backup@backdoor.com, system" },
    @{ Cell = "G30"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G31"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G32"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G33"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G37"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G38"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G39"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G40"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G41"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G42"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G44"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G45"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G46"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "H46"; Value = "39/57" },
    @{ Cell = "G56"; Value = "System, backup@backdoor.com, system" },
    @{ Cell = "G57"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G58"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G59"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G60"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G64"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G65"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G66"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G67"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G68"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G69"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G71"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G72"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G73"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "H73"; Value = "30/55" },
    @{ Cell = "G84"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G85"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G86"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G87"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G88"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G89"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G90"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G93"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G95"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G96"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G110"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G111"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G112"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G113"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G114"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G115"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G116"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G119"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G121"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G122"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G136"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G137"; Value = "System, backup@backdoor.com" },
    @{ Cell = "G138"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G139"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G140"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G141"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G142"; Value = "admin@admin.com, dnasr281@gmail.com" },
    @{ Cell = "G145"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G147"; Value = "System, dnasr281@gmail.com" },
    @{ Cell = "G148"; Value = "System, dnasr281@gmail.com" }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = $change.Value
}
